$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Changes 1 & 2: drop the spell-check "PMBoK" split (w:proofErr spellStart/
# spellEnd wrapping a separate run) by replacing the whole "(PMBoK)" span
# (which crosses the extra runs) with identical text. Find/Replace across
# run boundaries merges everything back into a single run and removes the
# now-irrelevant proofErr markers. wdReplaceAll (2) applies it to every
# occurrence in the document (there are exactly two, both needing the fix).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("(PMBoK)", $true, $false, $false, $false, $false, $true, 1, $false, "(PMBoK)", 2) | Out-Null

# ---------------------------------------------------------------------------
# Changes 3 & 4: highlight (bright green) the task bodies for items "3)" and
# "4)", and for item 3) split the leading "3) " marker so that only ") "
# (not the digit "3") picks up the highlight, matching item 4)'s style in
# the updated document. HighlightColorIndex only ever expands to cover
# whole runs it touches in this runtime, so the precise run layout is
# produced directly via Range.InsertXML with a flat-OPC wrapped fragment.
# ---------------------------------------------------------------------------

function Set-ParagraphRuns($paragraph, $bodyXml) {
    $start = $paragraph.Range.Start
    $end = $paragraph.Range.End
    # exclude the trailing paragraph mark so the existing pPr/paragraph
    # properties (spacing, numbering, ...) stay untouched
    $target = $d.Range($start, $end - 1)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $bodyXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($xml)
}

$p3 = $d.Paragraphs.Item(13)
$body3 = '<w:r><w:rPr><w:b/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:b/><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">) </w:t></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>Составьте Иерархическую структуру работ (ИСР, WBS) вашего проекта и её словарь. Обоснуйте выбранный подход к декомпозиции (продуктовый, функциональный, в соответствии с ЖЦ и т.д.).</w:t></w:r>'
Set-ParagraphRuns $p3 $body3

$p4 = $d.Paragraphs.Item(14)
$body4 = '<w:r><w:rPr><w:b/><w:highlight w:val="green"/></w:rPr><w:t>4)</w:t></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> Составьте реестр рисков для проекта. Укажите в нём оценки вероятности и влияния. Продумайте стратегии реагирования.</w:t></w:r>'
Set-ParagraphRuns $p4 $body4

Write-Host "edit complete"
